$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows {
    param($ws, $r1, $r2, $cols)
    foreach ($col in $cols) {
        $c1 = $ws.Cells.Item($r1, $col)
        $c2 = $ws.Cells.Item($r2, $col)
        $v1 = $c1.Value()
        $v2 = $c2.Value()
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

# Columns B=2, E=5, F=6, G=7 are swapped between the given row pairs
$cols = @(2, 5, 6, 7)

Swap-Rows $ws 183 184 $cols
Swap-Rows $ws 346 347 $cols
Swap-Rows $ws 351 352 $cols
Swap-Rows $ws 372 373 $cols
Swap-Rows $ws 379 380 $cols
Swap-Rows $ws 389 390 $cols
Swap-Rows $ws 419 420 $cols
Swap-Rows $ws 421 422 $cols
Swap-Rows $ws 579 580 $cols
Swap-Rows $ws 583 584 $cols
Swap-Rows $ws 586 587 $cols
Swap-Rows $ws 593 594 $cols
Swap-Rows $ws 601 602 $cols
Swap-Rows $ws 709 710 $cols
Swap-Rows $ws 715 716 $cols
Swap-Rows $ws 720 721 $cols
Swap-Rows $ws 872 873 $cols
